$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: fix the "118/9/14" text date in A6 into a real date value (11/9/2014)
$ws.Range("A6").Value = (New-Object DateTime(2014, 11, 9))

# Row 8: new timelog entry for Friday Sept 13th
$ws.Range("A8").Value = "13/9/14"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = "Writing down the design and architecture of the software as well as exploring the Twitter API"

# Copy formatting from row 7 (A7/B7/C7) down to row 8 so the new row matches the existing style
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

$ws.Range("C7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null

# Widen column C
$ws.Columns.Item(3).ColumnWidth = 70.1640625

# Update the active selection to C9
$ws.Range("C9").Select()
